$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2024-05-20T17:01:27+10:00"

# Update Publisher value (row 9, column B)
$ws.Range("B9").Value = "D Foulkes - Northern Australia Regional Digital Health Collaborative"

# Update Contact value (row 10, column B)
$ws.Range("B10").Value = "D Foulkes - Northern Australia Regional Digital Health Collaborative (https://nardhc.org)"

# Insert a new row at 11 for Jurisdiction / Australia, pushing Description/Purpose/Copyright/Immutable down
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "Australia"

# Copy the style from the row above (row 10) into new row 11 to match formatting
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Description (now row 12) gets the same text as Title
$ws.Range("B12").Value = "Medication Codes for Acute Rheumatic Fever"
